$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 6.896901285100976
$ws.Range("E2").Value = 6.139009826363375
$ws.Range("F2").Value = 6.576411474225445
$ws.Range("G2").Value = 5.772900851841581
$ws.Range("H2").Value = 6.874272700139254
$ws.Range("I2").Value = 6.124851285713226
$ws.Range("J2").Value = 6.500905559532678
$ws.Range("K2").Value = 5.639929978087127
